# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload"
#
# The FEB-22 sheet gets a new daily entry (row 19 / row r="27"-"28" in the
# sheet XML): date 23-Feb-2022, application "RPA GSS", two wrapped comment
# lines with their own % complete / status, and the view scrolls down to
# where the new row was typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet          # "FEB-22" is already the active sheet

# Row 27 (B column) needs the same date formatting (style) as the rows
# above it. Copy just the format from B25 (a date cell) onto B27 instead of
# assigning a NumberFormat string, so it reuses the workbook's existing
# date style rather than minting a new one.
$ws.Range("B25").Copy()
$ws.Range("B27").PasteSpecial(-4122)   # xlPasteFormats

# New entry: No.19, 23-Feb-2022, RPA GSS
$ws.Range("A27").Value = 19
$ws.Range("B27").Value = 44615
$ws.Range("C27").Value = "RPA GSS"
$ws.Range("D27").Value = "1. Updating master data file  is work in progress as new suggestion has been received from Mohan san, ( we have to extract data pagewise and update to be done based on the configuration file)"
$ws.Range("E27").Value = 0.75
$ws.Range("F27").Value = "WIP"

# Second comment line for the same entry, continuing on row 28
$ws.Range("D28").Value = "2. Supported to Captcha issue at token system  (RPA GSS)( Now the login is getting success from the six to ten captcha images) and still the task to be tested."
$ws.Range("E28").Value = 0.6
$ws.Range("F28").Value = "WIP"

# Both rows grew to a two-line wrap height once the long comments were typed in
$ws.Rows.Item(27).RowHeight = 28.8
$ws.Rows.Item(28).RowHeight = 28.8

# Scroll/selection state left behind after typing the new row
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D34").Select()
